$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main coils")

$ws.Range("C2").Value = 'NE CORNER OFFICE SZ-CAV UNITARY'
$ws.Range("C3").Value = 'NW CORNER OFFICE SZ-CAV UNITARY'
$ws.Range("C4").Value = 'SE CORNER OFFICE SZ-CAV UNITARY'
$ws.Range("C5").Value = 'SW CORNER OFFICE SZ-CAV UNITARY'
$ws.Range("C6").Value = 'CORE AUDITORIUM SZ-CAV UNITARY'
$ws.Range("C7").Value = 'SHOP EL1 WEST PERIM SPC (G.W1) SZ-CAV UNITARY'
$ws.Range("C8").Value = 'CLASSROOM EL2 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C9").Value = 'CLASSROOM EL2 SOUTH PERIM SPC (G.S2) SZ-CAV UNITARY'
$ws.Range("C10").Value = 'CLASSROOM EL2 WEST PERIM SPC (G.W1) SZ-CAV UNITARY'
$ws.Range("C11").Value = 'DINING EL3 SW PERIM SPC (G.SW1) SZ-CAV UNITARY'
$ws.Range("C12").Value = 'CLASSROOM EL4 EAST PERIM SPC (G.E2) SZ-CAV UNITARY'
$ws.Range("C13").Value = 'CLASSROOM EL4 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C14").Value = 'CLASSROOM EL4 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C15").Value = 'CLASSROOM EL4 WEST PERIM SPC (G.W4) SZ-CAV UNITARY'
$ws.Range("C16").Value = 'OFFICE EL4 CORE SPC (G.C6) SZ-CAV UNITARY'
$ws.Range("C17").Value = 'CLASSROOM EL4 EAST PERIM SPC (T.E14) SZ-CAV UNITARY'
$ws.Range("C18").Value = 'CLASSROOM EL4 NORTH PERIM SPC (T.N15) SZ-CAV UNITARY'
$ws.Range("C19").Value = 'CLASSROOM EL4 SOUTH PERIM SPC (T.S13) SZ-CAV UNITARY'
$ws.Range("C20").Value = 'CLASSROOM EL4 WEST PERIM SPC (T.W16) SZ-CAV UNITARY'
$ws.Range("C21").Value = 'OFFICE EL4 CORE SPC (T.C18) SZ-CAV UNITARY'
$ws.Range("C22").Value = 'CLASSROOM G.E1 SZ-CAV UNITARY'
$ws.Range("C23").Value = 'CLASSROOM G.NNE2 SZ-CAV UNITARY'
$ws.Range("C24").Value = 'CLASSROOM G.SSE3 SZ-CAV UNITARY'
$ws.Range("C25").Value = 'CLASSROOM G.W4 SZ-CAV UNITARY'
$ws.Range("C26").Value = 'DINING G.SSW1 SZ-CAV UNITARY'
$ws.Range("C27").Value = 'GYMNASIUM G.N2 SZ-CAV UNITARY'
$ws.Range("C28").Value = 'CLASSROOM EL1 SPC (G.1) SZ-CAV UNITARY'
$ws.Range("C29").Value = 'CLASSROOM E2 NORTH PERIM (G.N3) SZ-CAV UNITARY'
$ws.Range("C30").Value = 'CLASSROOM E2 SOUTH PERIM (G.S2) SZ-CAV UNITARY'
$ws.Range("C31").Value = 'CLASSROOM E2 WEST PERIM (G.W1) SZ-CAV UNITARY'
$ws.Range("C32").Value = 'CLASSROOM E4 CORE (G.C6) SZ-CAV UNITARY'
$ws.Range("C33").Value = 'CLASSROOM E4 EAST PERIM (G.E2) SZ-CAV UNITARY'
$ws.Range("C34").Value = 'CLASSROOM E4 NORTH PERIM (G.N3) SZ-CAV UNITARY'
$ws.Range("C35").Value = 'CLASSROOM E4 SOUTH PERIM (G.S1) SZ-CAV UNITARY'
$ws.Range("C36").Value = 'CLASSROOM E4 WEST PERIM (G.W4) SZ-CAV UNITARY'
$ws.Range("C37").Value = 'CLASSROOM E2 NORTH PERIM (G.N3) SZ-CAV UNITARY'
$ws.Range("C38").Value = 'CLASSROOM E2 SOUTH PERIM (G.S2) SZ-CAV UNITARY'
$ws.Range("C39").Value = 'CLASSROOM E2 WEST PERIM (G.W1) SZ-CAV UNITARY'
$ws.Range("C40").Value = 'DINING E3 SOUTH PERIM (G.S1) SZ-CAV UNITARY'
$ws.Range("C41").Value = 'MAINTENANCE E1 WEST PERIM (G.W1) SZ-CAV UNITARY'
$ws.Range("C42").Value = 'CLASSROOM E4 EAST PERIM (G.E2) SZ-CAV UNITARY'
$ws.Range("C43").Value = 'CLASSROOM E4 NORTH PERIM (G.N3) SZ-CAV UNITARY'
$ws.Range("C44").Value = 'CLASSROOM E4 SOUTH PERIM (G.S1) SZ-CAV UNITARY'
$ws.Range("C45").Value = 'CLASSROOM E4 WEST PERIM (G.W4) SZ-CAV UNITARY'
$ws.Range("C46").Value = 'OFFICE E4 CORE (G.C6) SZ-CAV UNITARY'
$ws.Range("C47").Value = 'CLASSROOM E4 EAST PERIM (M.E14) SZ-CAV UNITARY'
$ws.Range("C48").Value = 'CLASSROOM E4 NORTH PERIM (M.N15) SZ-CAV UNITARY'
$ws.Range("C49").Value = 'CLASSROOM E4 SOUTH PERIM (M.S13) SZ-CAV UNITARY'
$ws.Range("C50").Value = 'CLASSROOM E4 WEST PERIM (M.W16) SZ-CAV UNITARY'
$ws.Range("C51").Value = 'OFFICE E4 CORE (M.C18) SZ-CAV UNITARY'
$ws.Range("C52").Value = 'CLASSROOM E4 EAST PERIM (T.E26) SZ-CAV UNITARY'
$ws.Range("C53").Value = 'CLASSROOM E4 NORTH PERIM (T.N27) SZ-CAV UNITARY'
$ws.Range("C54").Value = 'CLASSROOM E4 SOUTH PERIM (T.S25) SZ-CAV UNITARY'
$ws.Range("C55").Value = 'CLASSROOM E4 WEST PERIM (T.W28) SZ-CAV UNITARY'
$ws.Range("C56").Value = 'OFFICE E4 CORE (T.C30) SZ-CAV UNITARY'
$ws.Range("C57").Value = 'CORRIDOR E9 CORE (G.C3) SZ-CAV UNITARY'
$ws.Range("C58").Value = 'CORRIDOR E9 CORE (M.C6) SZ-CAV UNITARY'
$ws.Range("C59").Value = 'CORRIDOR E9 CORE (T.C9) SZ-CAV UNITARY'
$ws.Range("C60").Value = 'GROCERY SALES EL1 SSE PERIM SPC (G.SSE1) SZ-CAV UNITARY'
$ws.Range("C61").Value = 'INDLOADDOCK EL1 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C62").Value = 'OFFICE EL1 SW PERIM SPC (G.SW2) SZ-CAV UNITARY'
$ws.Range("C63").Value = 'REFFOODPREP EL1 CORE SPC (G.C4) SZ-CAV UNITARY'
$ws.Range("C64").Value = 'REFWALKINCOOL EL1 CORE SPC (G.C6) SZ-CAV UNITARY'
$ws.Range("C65").Value = 'REFWALKINFREEZE EL1 CORE SPC (G.C5) SZ-CAV UNITARY'
$ws.Range("C66").Value = 'DINING EL1 NORTH PERIM (G.N1) SZ-CAV UNITARY'
$ws.Range("C67").Value = 'DINING EL1 SOUTH PERIM (G.S2) SZ-CAV UNITARY'
$ws.Range("C68").Value = 'OFFICE EL2 CORE (G.C7) SZ-CAV UNITARY'
$ws.Range("C69").Value = 'OFFICE EL2 EAST PERIM (G.E1) SZ-CAV UNITARY'
$ws.Range("C70").Value = 'OFFICE EL2 NORTH PERIM (G.N6) SZ-CAV UNITARY'
$ws.Range("C71").Value = 'OFFICE EL2 SOUTH PERIM (G.S5) SZ-CAV UNITARY'
$ws.Range("C72").Value = 'SURGERY EL2 NNW PERIM (G.NNW3) SZ-CAV UNITARY'
$ws.Range("C73").Value = 'SURGERY EL2 SSW PERIM (G.SSW2) SZ-CAV UNITARY'
$ws.Range("C74").Value = 'SURGERY EL2 WEST PERIM (G.W4) SZ-CAV UNITARY'
$ws.Range("C75").Value = 'BARCASINO EL2 NORTH PERIM SPC (G.N2) SZ-CAV UNITARY'
$ws.Range("C76").Value = 'BARCASINO EL3 NORTH PERIM SPC (G.N1) SZ-CAV UNITARY'
$ws.Range("C77").Value = 'CORRIDOR EL4 CORE SPC (G.C5) SZ-CAV UNITARY'
$ws.Range("C78").Value = 'CORRIDOR EL4 CORE SPC (M.C15) SZ-CAV UNITARY'
$ws.Range("C79").Value = 'CORRIDOR EL4 CORE SPC (T.C25) SZ-CAV UNITARY'
$ws.Range("C80").Value = 'DINING EL2 WEST PERIM SPC (G.W3) SZ-CAV UNITARY'
$ws.Range("C81").Value = 'DINING EL3 WEST PERIM SPC (G.W3) SZ-CAV UNITARY'
$ws.Range("C82").Value = 'LAUNDRY EL1 WEST PERIM SPC (G.W1) SZ-CAV UNITARY'
$ws.Range("C83").Value = 'LOBBY EL1 WEST PERIM SPC (G.W2) SZ-CAV UNITARY'
$ws.Range("C84").Value = 'OFFICE EL1 EAST PERIM SPC (G.E3) SZ-CAV UNITARY'
$ws.Range("C85").Value = 'CONFERENCE EL1 CORE SPC (G.C7) SZ-CAV UNITARY'
$ws.Range("C86").Value = 'CORRIDOR EL1 CORE SPC (G.C6) SZ-CAV UNITARY'
$ws.Range("C87").Value = 'DINING EL1 CORE SPC (G.C9) SZ-CAV UNITARY'
$ws.Range("C88").Value = 'OFFICE EL1 CORE SPC (G.C11) SZ-CAV UNITARY'
$ws.Range("C89").Value = 'OFFICE EL1 EAST PERIM SPC (G.E2) SZ-CAV UNITARY'
$ws.Range("C90").Value = 'OFFICE EL1 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C91").Value = 'OFFICE EL1 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C92").Value = 'OFFICE EL1 WEST PERIM SPC (G.W4) SZ-CAV UNITARY'
$ws.Range("C93").Value = 'STOCKROOM EL1 NE PERIM SPC (G.NE2) SZ-CAV UNITARY'
$ws.Range("C94").Value = 'STOCKROOM EL1 SW PERIM SPC (G.SW3) SZ-CAV UNITARY'
$ws.Range("C95").Value = 'WORK EL1 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C96").Value = 'CORRIDOR EL1 CORE PERIM SPC (G.C5) SZ-CAV UNITARY'
$ws.Range("C97").Value = 'CORRIDOR EL1 SOUTH PERIM SPC (G.S9) SZ-CAV UNITARY'
$ws.Range("C98").Value = 'CORRIDOR EL2 CORE SPC (G.C5) SZ-CAV UNITARY'
$ws.Range("C99").Value = 'CORRIDOR EL2 NNE PERIM SPC (G.NNE4) SZ-CAV UNITARY'
$ws.Range("C100").Value = 'CORRIDOR EL2 NORTH PERIM SPC (G.N1) SZ-CAV UNITARY'
$ws.Range("C101").Value = 'CORRIDOR EL2 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C102").Value = 'CORRIDOR EL2 SSW PERIM SPC (G.SSW2) SZ-CAV UNITARY'
$ws.Range("C103").Value = 'CORRIDOR EL2 SOUTH PERIM SPC (G.S6) SZ-CAV UNITARY'
$ws.Range("C104").Value = 'CORRIDOR EL2 SOUTH PERIM SPC (G.S7) SZ-CAV UNITARY'
$ws.Range("C105").Value = 'LAUNDRY EL1 SOUTH PERIM SPC (G.S7) SZ-CAV UNITARY'
$ws.Range("C106").Value = 'OFFICE EL1 SE PERIM SPC (G.SE6) SZ-CAV UNITARY'
$ws.Range("C107").Value = 'DINING EL1 WEST PERIM SPC (G.W1) SZ-CAV UNITARY'
$ws.Range("C108").Value = 'OFFICE EL2 CORE SPC (G.C1) SZ-CAV UNITARY'
$ws.Range("C109").Value = 'OFFICE EL2 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C110").Value = 'OFFICE EL2 SOUTH PERIM SPC (G.S2) SZ-CAV UNITARY'
$ws.Range("C111").Value = 'OFFICE EL3 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C112").Value = 'OFFICE EL3 SOUTH PERIM SPC (G.S2) SZ-CAV UNITARY'
$ws.Range("C113").Value = 'OFFICE EL3 WEST PERIM SPC (G.W1) SZ-CAV UNITARY'
$ws.Range("C114").Value = 'OFFICEOPEN EL1 CORE SPC (M.C15) SZ-CAV UNITARY'
$ws.Range("C115").Value = 'OFFICEOPEN EL1 CORE SPC (T.C25) SZ-CAV UNITARY'
$ws.Range("C116").Value = 'OFFICESMALL EL1 EAST PERIM SPC (M.E12) SZ-CAV UNITARY'
$ws.Range("C117").Value = 'OFFICESMALL EL1 EAST PERIM SPC (T.E22) SZ-CAV UNITARY'
$ws.Range("C118").Value = 'OFFICESMALL EL1 NORTH PERIM SPC (M.N13) SZ-CAV UNITARY'
$ws.Range("C119").Value = 'OFFICESMALL EL1 NORTH PERIM SPC (T.N23) SZ-CAV UNITARY'
$ws.Range("C120").Value = 'OFFICESMALL EL1 SOUTH PERIM SPC (M.S11) SZ-CAV UNITARY'
$ws.Range("C121").Value = 'OFFICESMALL EL1 SOUTH PERIM SPC (T.S21) SZ-CAV UNITARY'
$ws.Range("C122").Value = 'OFFICESMALL EL1 WEST PERIM SPC (M.W14) SZ-CAV UNITARY'
$ws.Range("C123").Value = 'OFFICESMALL EL1 WEST PERIM SPC (T.W24) SZ-CAV UNITARY'
$ws.Range("C124").Value = 'HALL EL1 CORE SPC (G.C5) SZ-CAV UNITARY'
$ws.Range("C125").Value = 'OFFICESMALL EL1 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C126").Value = 'OFFICESMALL EL1 EAST PERIM SPC (G.E2) SZ-CAV UNITARY'
$ws.Range("C127").Value = 'OFFICESMALL EL1 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C128").Value = 'OFFICESMALL EL1 WEST PERIM SPC (G.W4) SZ-CAV UNITARY'
$ws.Range("C129").Value = 'HALL EL1 CORE SPC (T.C15) SZ-CAV UNITARY'
$ws.Range("C130").Value = 'OFFICESMALL EL1 EAST PERIM SPC (T.E12) SZ-CAV UNITARY'
$ws.Range("C131").Value = 'OFFICESMALL EL1 NORTH PERIM SPC (T.N13) SZ-CAV UNITARY'
$ws.Range("C132").Value = 'OFFICESMALL EL1 SOUTH PERIM SPC (T.S11) SZ-CAV UNITARY'
$ws.Range("C133").Value = 'OFFICESMALL EL1 WEST PERIM SPC (T.W14) SZ-CAV UNITARY'
$ws.Range("C134").Value = 'DINING EL1 ESE PERIM SPC (G.ESE1) SZ-CAV UNITARY'
$ws.Range("C135").Value = 'LOBBYWAITING EL1 SSW PERIM SPC (G.SSW2) SZ-CAV UNITARY'
$ws.Range("C136").Value = 'RESTROOM EL1 NORTH PERIM SPC (G.N4) SZ-CAV UNITARY'
$ws.Range("C137").Value = 'LOBBYWAITING EL1 SW PERIM SPC (G.SW2) SZ-CAV UNITARY'
$ws.Range("C138").Value = 'RESTROOM EL1 NORTH PERIM SPC (G.N4) SZ-CAV UNITARY'
$ws.Range("C139").Value = 'EL1 CORE SPC (G.C5) SZ-CAV UNITARY'
$ws.Range("C140").Value = 'RETAILSALES EL1 NORTH PERIM SPC (G.N3) SZ-CAV UNITARY'
$ws.Range("C141").Value = 'RETAILSALES EL1 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C142").Value = 'RETAILSALES EL1 WEST PERIM SPC (G.W4) SZ-CAV UNITARY'
$ws.Range("C143").Value = 'RETAILSALES EL1 CORE SPC (M.C15) SZ-CAV UNITARY'
$ws.Range("C144").Value = 'RETAILSALES EL1 EAST PERIM SPC (M.E12) SZ-CAV UNITARY'
$ws.Range("C145").Value = 'RETAILSALES EL1 NORTH PERIM SPC (M.N13) SZ-CAV UNITARY'
$ws.Range("C146").Value = 'RETAILSALES EL1 SOUTH PERIM SPC (M.S11) SZ-CAV UNITARY'
$ws.Range("C147").Value = 'RETAILSALES EL1 WEST PERIM SPC (M.W14) SZ-CAV UNITARY'
$ws.Range("C148").Value = 'RETAILSALES EL1 CORE SPC (T.C25) SZ-CAV UNITARY'
$ws.Range("C149").Value = 'RETAILSALES EL1 EAST PERIM SPC (T.E22) SZ-CAV UNITARY'
$ws.Range("C150").Value = 'RETAILSALES EL1 NORTH PERIM SPC (T.N23) SZ-CAV UNITARY'
$ws.Range("C151").Value = 'RETAILSALES EL1 SOUTH PERIM SPC (T.S21) SZ-CAV UNITARY'
$ws.Range("C152").Value = 'RETAILSALES EL1 WEST PERIM SPC (T.W24) SZ-CAV UNITARY'
$ws.Range("C153").Value = 'AUTOSHOP EL1 WEST PERIM SPC (G.W5) SZ-CAV UNITARY'
$ws.Range("C154").Value = 'OFFICE EL1 EAST PERIM SPC (G.E4) SZ-CAV UNITARY'
$ws.Range("C155").Value = 'RETAIL EL1 CORE SPC (G.C8) SZ-CAV UNITARY'
$ws.Range("C156").Value = 'RETAIL EL1 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C157").Value = 'RETAIL EL1 WSW PERIM SPC (G.WSW7) SZ-CAV UNITARY'
$ws.Range("C158").Value = 'STOCK EL1 EAST PERIM SPC (G.E3) SZ-CAV UNITARY'
$ws.Range("C159").Value = 'STOCK EL1 NORTH PERIM SPC (G.N2) SZ-CAV UNITARY'
$ws.Range("C160").Value = 'RETAILSALES EL1 SOUTH PERIM (G.S1) SZ-CAV UNITARY'
$ws.Range("C161").Value = 'STOCKROOM EL1 NORTH PERIM (G.N2) SZ-CAV UNITARY'
$ws.Range("C162").Value = 'WAREHOUSECOND EL1 NORTH PERIM SPC (G.N2) SZ-CAV UNITARY'
$ws.Range("C163").Value = 'WAREHOUSECOND EL1 SOUTH PERIM SPC (G.S1) SZ-CAV UNITARY'
$ws.Range("C1").Value = 'object name'
$ws.Range("A2").Select()
